$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.152.48'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '1.859.20'
$ws.Range("E3").Value = '  -0.56%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.56'
$ws.Range("E5").Value = '  -0.77%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4679'
$ws.Range("E7").Value = '  -0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.90'
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2834'
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06465'
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.02'
$ws.Range("E11").Value = '  -3.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07734'
$ws.Range("E12").Value = '  -3.55%  '
$ws.Range("D13").Value = '1.870.90'
$ws.Range("E13").Value = '  +0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.59'
$ws.Range("E14").Value = '  -3.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.055'
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.6793'
$ws.Range("E16").Value = '  -0.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '266.51'
$ws.Range("E17").Value = '  -0.90%  '
$ws.Range("D18").Value = '30.113.61'
$ws.Range("E18").Value = '  -0.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.32'
$ws.Range("E19").Value = '  -4.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007577'
$ws.Range("E20").Value = '  -1.37%  '
$ws.Range("E21").Value = '  -0.08%  '
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.147'
$ws.Range("E23").Value = '  -2.35%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.100'
$ws.Range("E24").Value = '  -1.85%  '
$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.288'
$ws.Range("E25").Value = '  -1.54%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.25'
$ws.Range("E26").Value = '  -1.98%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.52'
$ws.Range("E27").Value = '  -2.05%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.884'
$ws.Range("E28").Value = '  -3.35%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.365'
$ws.Range("E29").Value = '  -0.34%  '
$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09841'
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.452'
$ws.Range("E31").Value = '  -0.87%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.203'
$ws.Range("E32").Value = '  -4.05%  '
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.976'
$ws.Range("E33").Value = '  -2.36%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04652'
$ws.Range("E34").Value = '  -1.05%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.114'
$ws.Range("E35").Value = '  -1.91%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6864'
$ws.Range("E36").Value = '  -2.05%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.713'
$ws.Range("E37").Value = '  +0.28%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01813'
$ws.Range("E38").Value = '  -3.45%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.717'
$ws.Range("E39").Value = '  +3.69%  '
$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.268'
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("B41").Value = 'Aave'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '70.66'
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9998'
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8322'
$ws.Range("E43").Value = '  -1.15%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.878'
$ws.Range("E44").Value = '  -4.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '101.92'
$ws.Range("E45").Value = '  -1.20%  '
$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4044'
$ws.Range("E46").Value = '  -3.02%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.134'
$ws.Range("E47").Value = '  -0.71%  '
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '922.84'
$ws.Range("E48").Value = '  +0.38%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.940'
$ws.Range("E49").Value = '  -1.89%  '
$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.07'
$ws.Range("E50").Value = '  -1.19%  '
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05557'
$ws.Range("E51").Value = '  -2.20%  '
